$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "jYFrx492"
$ws.Range("B2").Value = 23072635
$ws.Range("C2").Value = "kmnmcik32"
$ws.Range("D2").Value = "c9J!Sz#6"
$ws.Range("F2").Value = "PjcZTUGf"
$ws.Range("G2").Value = "gair"
